$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 16.797
$ws.Range("D8").Value = -8.559000000000001
$ws.Range("D10").Value = -8.231999999999999
$ws.Range("D12").Value = -7.289999999999999
$ws.Range("E14").Value = 17.035
$ws.Range("E15").Value = 16.367
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("E18").Value = 16.353
$ws.Range("E20").Value = 16.459
$ws.Range("D25").Value = -8.228999999999999
$ws.Range("E29").Value = 17.05
$ws.Range("E30").Value = 16.474
$ws.Range("E31").Value = 16.527
$ws.Range("E35").Value = 16.552
$ws.Range("D37").Value = -8.175000000000001
$ws.Range("E40").Value = 16.627
$ws.Range("E44").Value = 16.815
$ws.Range("E50").Value = 16.319
$ws.Range("E54").Value = 16.692
$ws.Range("D55").Value = -8.318999999999999
$ws.Range("D68").Value = -7.188
$ws.Range("E68").Value = 17.777
$ws.Range("E76").Value = 16.558
$ws.Range("D77").Value = -7.936
$ws.Range("D78").Value = -7.983999999999999
$ws.Range("D79").Value = -7.609
$ws.Range("D80").Value = -7.923
$ws.Range("D81").Value = -7.536
$ws.Range("D82").Value = -8.154
$ws.Range("D84").Value = -8.504999999999999
$ws.Range("E87").Value = 16.366
$ws.Range("E88").Value = 16.278
$ws.Range("E92").Value = 17.658
$ws.Range("E96").Value = 16.325
$ws.Range("E98").Value = 16.299
$ws.Range("D101").Value = -7.845000000000001
$ws.Range("E101").Value = 16.54
$ws.Range("D102").Value = -8.098000000000001
$ws.Range("E102").Value = 16.674
